$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3123.25
$ws.Range("I58").Value = 80.666664
$ws.Range("J58").Value = 3438
$ws.Range("K58").Value = 241.999992
$ws.Range("L58").Value = 10314
$ws.Range("M58").Value = -91.99999199999999
$ws.Range("N58").Value = -10614
$ws.Range("H70").Value = 1557.05
$ws.Range("I70").Value = 1900.2858
$ws.Range("J70").Value = 1372.2307
$ws.Range("K70").Value = 5700.857400000001
$ws.Range("L70").Value = 4116.6921
$ws.Range("M70").Value = -5430.857400000001
$ws.Range("N70").Value = -4656.6921
$ws.Range("H73").Value = 1557.05
$ws.Range("I73").Value = 1900.2858
$ws.Range("J73").Value = 1372.2307
$ws.Range("K73").Value = 5700.857400000001
$ws.Range("L73").Value = 4116.6921
$ws.Range("M73").Value = -4764.857400000001
$ws.Range("N73").Value = -5988.6921
$ws.Range("H74").Value = 2276033.2
$ws.Range("I74").Value = 2705997.2
$ws.Range("J74").Value = 3367.1428
$ws.Range("K74").Value = 2705997.2
$ws.Range("L74").Value = 3367.1428
$ws.Range("M74").Value = -2705061.2
$ws.Range("N74").Value = -5239.1428
$ws.Range("H77").Value = 2276033.2
$ws.Range("I77").Value = 2705997.2
$ws.Range("J77").Value = 3367.1428
$ws.Range("K77").Value = 13529986
$ws.Range("L77").Value = 16835.714
$ws.Range("M77").Value = -13525306
$ws.Range("N77").Value = -26195.714
$ws.Range("H129").Value = 2253.4
$ws.Range("I129").Value = 612.5714
$ws.Range("J129").Value = 2752.7827
$ws.Range("K129").Value = 1837.7142
$ws.Range("L129").Value = 8258.348100000001
$ws.Range("M129").Value = 3162.2858
$ws.Range("N129").Value = -18258.3481
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1593.0416
$ws.Range("I45").Value = 1496.3182
$ws.Range("J45").Value = 2657
$ws.Range("K45").Value = 1496.3182
$ws.Range("L45").Value = 2657
$ws.Range("M45").Value = -1119.3182
$ws.Range("N45").Value = -3411
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7042.2
$ws.Range("I94").Value = 12440.444
$ws.Range("J94").Value = 2625.4546
$ws.Range("K94").Value = 12440.444
$ws.Range("L94").Value = 2625.4546
$ws.Range("M94").Value = -11989.444
$ws.Range("N94").Value = -3527.4546
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 8759.286
$ws.Range("I41").Value = 712.5
$ws.Range("J41").Value = 19488.334
$ws.Range("K41").Value = 712.5
$ws.Range("L41").Value = 19488.334
$ws.Range("M41").Value = -284.5
$ws.Range("N41").Value = -20344.334
$ws.Range("H59").Value = 10300
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 10300
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 10300
$ws.Range("N59").Value = -12590
$ws.Range("H60").Value = 7673.25
$ws.Range("I60").Value = 4397.6665
$ws.Range("J60").Value = 17500
$ws.Range("K60").Value = 4397.6665
$ws.Range("L60").Value = 17500
$ws.Range("M60").Value = -3886.6665
$ws.Range("N60").Value = -18522
$ws.Range("H68").Value = 26173.2
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 26173.2
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 26173.2
$ws.Range("N68").Value = -27671.2
$ws.Range("H71").Value = 26173.2
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 26173.2
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 78519.60000000001
$ws.Range("N71").Value = -86007.60000000001
$ws.Range("H74").Value = 10900
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10900
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 10900
$ws.Range("N74").Value = -12648
$ws.Range("H77").Value = 10900
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10900
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 32700
$ws.Range("N77").Value = -41436
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2678.3333
$ws.Range("I63").Value = 2007
$ws.Range("J63").Value = 3014
$ws.Range("K63").Value = 6021
$ws.Range("L63").Value = 9042
$ws.Range("M63").Value = -5272
$ws.Range("N63").Value = -10540
$ws.Range("H64").Value = 2759159
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 3034974.8
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 9104924.399999999
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -9105464.399999999
$ws.Range("H66").Value = 2678.3333
$ws.Range("I66").Value = 2007
$ws.Range("J66").Value = 3014
$ws.Range("K66").Value = 18063
$ws.Range("L66").Value = 27126
$ws.Range("M66").Value = -14319
$ws.Range("N66").Value = -34614
$ws.Range("H67").Value = 2759159
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 3034974.8
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 9104924.399999999
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -9106796.399999999
$ws.Range("H68").Value = 2327.8975
$ws.Range("I68").Value = 4461.4
$ws.Range("J68").Value = 1321.5283
$ws.Range("K68").Value = 13384.2
$ws.Range("L68").Value = 3964.5849
$ws.Range("M68").Value = -12573.2
$ws.Range("N68").Value = -5586.5849
$ws.Range("H70").Value = 2529.7334
$ws.Range("I70").Value = 789.6
$ws.Range("J70").Value = 3399.8
$ws.Range("K70").Value = 2368.8
$ws.Range("L70").Value = 10199.4
$ws.Range("M70").Value = -2053.8
$ws.Range("N70").Value = -10829.4
$ws.Range("H71").Value = 2327.8975
$ws.Range("I71").Value = 4461.4
$ws.Range("J71").Value = 1321.5283
$ws.Range("K71").Value = 40152.6
$ws.Range("L71").Value = 11893.7547
$ws.Range("M71").Value = -36096.6
$ws.Range("N71").Value = -20005.7547
$ws.Range("H73").Value = 2529.7334
$ws.Range("I73").Value = 789.6
$ws.Range("J73").Value = 3399.8
$ws.Range("K73").Value = 2368.8
$ws.Range("L73").Value = 10199.4
$ws.Range("M73").Value = -1276.8
$ws.Range("N73").Value = -12383.4
$ws.Range("H87").Value = 338672
$ws.Range("I87").Value = 1000
$ws.Range("J87").Value = 507508
$ws.Range("K87").Value = 3000
$ws.Range("L87").Value = 1522524
$ws.Range("M87").Value = -1752
$ws.Range("N87").Value = -1525020
$ws.Range("H90").Value = 338672
$ws.Range("I90").Value = 1000
$ws.Range("J90").Value = 507508
$ws.Range("K90").Value = 9000
$ws.Range("L90").Value = 4567572
$ws.Range("M90").Value = -2760
$ws.Range("N90").Value = -4580052
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1017.64703
$ws.Range("I81").Value = 1006.25
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 2012.5
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -951.5
$ws.Range("N81").Value = -4522
$ws.Range("H84").Value = 1017.64703
$ws.Range("I84").Value = 1006.25
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 10062.5
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -4758.5
$ws.Range("N84").Value = -22608
$ws.Range("H107").Value = 1186.6111
$ws.Range("I107").Value = 1047.1666
$ws.Range("J107").Value = 1465.5
$ws.Range("K107").Value = 3141.4998
$ws.Range("L107").Value = 4396.5
$ws.Range("M107").Value = -1221.4998
$ws.Range("N107").Value = -8236.5
